# Weekly update: insert a new price record for "Vega Modelo de Temuco -
# Espárragos" at row 100 (pushing the existing rows 100-118 down to
# 101-119), matching the "Fruta / hortaliza, semanal" refresh pattern.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 100..118 down to 101..119.
$ws.Rows.Item(100).Insert()

# Populate the newly inserted row 100 with the latest weekly record.
$ws.Cells.Item(100, 1).Value  = 10
$ws.Cells.Item(100, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(100, 3).Value  = "La Araucanía"
$ws.Cells.Item(100, 4).Value  = 45244
$ws.Cells.Item(100, 5).Value  = 9
$ws.Cells.Item(100, 6).Value  = 300000000
$ws.Cells.Item(100, 7).Value  = "Espárragos"
$ws.Cells.Item(100, 8).Value  = "Sin especificar"
$ws.Cells.Item(100, 9).Value  = "Primera"
$ws.Cells.Item(100, 10).Value = 250
$ws.Cells.Item(100, 11).Value = 1600
$ws.Cells.Item(100, 12).Value = 1600
$ws.Cells.Item(100, 13).Value = 1600
$ws.Cells.Item(100, 14).Value = "$/kilo"
$ws.Cells.Item(100, 15).Value = "Región del Maule"
$ws.Cells.Item(100, 16).Value = 1600
$ws.Cells.Item(100, 17).Value = 1
$ws.Cells.Item(100, 18).Value = "Hortaliza"
